# Applies updated Fitness (column C) values for rows 2-54 on Sheet1,
# matching the target diff for run_3.xlsx (mp0-8_4 log).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fitnessUpdates = @{
    2 = 12362
    3 = 12362
    4 = 12012
    5 = 12012
    6 = 10617
    7 = 10617
    8 = 10617
    9 = 10617
    10 = 10617
    11 = 10617
    12 = 10472
    13 = 10472
    14 = 9940
    15 = 9576
    16 = 9576
    17 = 9576
    18 = 9452
    19 = 9452
    20 = 9452
    21 = 9452
    22 = 9452
    23 = 9452
    24 = 9452
    25 = 9452
    26 = 9041
    27 = 9041
    28 = 9041
    29 = 9041
    30 = 9041
    31 = 9041
    32 = 8844
    33 = 8844
    34 = 8844
    35 = 8844
    36 = 8844
    37 = 8844
    38 = 8844
    39 = 8710
    40 = 8710
    41 = 8710
    42 = 8710
    43 = 8710
    44 = 8710
    45 = 8710
    46 = 8372
    47 = 8372
    48 = 8372
    49 = 8372
    50 = 8372
    51 = 7917
    52 = 7917
    53 = 7917
    54 = 7917
}

foreach ($row in $fitnessUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $fitnessUpdates[$row]
}
